$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.814.10'
$ws.Range("E2").Value = '  -0.02%  '
$ws.Range("D3").Value = '3.812.62'
$ws.Range("E3").Value = '  +1.77%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '615.47'
$ws.Range("E5").Value = '  -0.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '177.71'
$ws.Range("E6").Value = '  -0.99%  '
$ws.Range("D7").Value = '3.810.45'
$ws.Range("E7").Value = '  +1.78%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -1.13%  '
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.49'
$ws.Range("E11").Value = '  +2.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.486'
$ws.Range("E12").Value = '  -0.58%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '39.98'
$ws.Range("E13").Value = '  -2.53%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000256'
$ws.Range("E14").Value = '  -1.76%  '
$ws.Range("D15").Value = '4.442.05'
$ws.Range("E15").Value = '  +1.98%  '
$ws.Range("D16").Value = '3.810.39'
$ws.Range("E16").Value = '  +1.79%  '
$ws.Range("D17").Value = '69.801.23'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.59'
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("E19").Value = '  -3.57%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '16.71'
$ws.Range("E20").Value = '  +0.29%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '511.60'
$ws.Range("E21").Value = '  +1.10%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.64'
$ws.Range("E22").Value = '  +2.11%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.737'
$ws.Range("E23").Value = '  +1.13%  '
$ws.Range("E24").Value = '  -1.18%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.57'
$ws.Range("E25").Value = '  -0.62%  '
$ws.Range("E26").Value = '  +7.43%  '
$ws.Range("E27").Value = '  -2.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.59'
$ws.Range("E28").Value = '  -4.81%  '
$ws.Range("E29").Value = '  +0.04%  '
$ws.Range("E30").Value = '  +0.70%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.99'
$ws.Range("E31").Value = '  +2.72%  '
$ws.Range("E32").Value = '  +2.52%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '31.69'
$ws.Range("E33").Value = '  +1.32%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.114'
$ws.Range("E34").Value = '  -1.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("E35").Value = '  +0.13%  '
$ws.Range("E36").Value = '  -1.14%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.15'
$ws.Range("E37").Value = '  -0.81%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.141'
$ws.Range("E38").Value = '  +6.81%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '480.05'
$ws.Range("E39").Value = '  +13.82%  '
$ws.Range("E40").Value = '  +1.45%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.07'
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("E42").Value = '  +5.14%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '49.79'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '44.16'
$ws.Range("E44").Value = '  -3.01%  '
$ws.Range("E45").Value = '  -1.21%  '
$ws.Range("D46").Value = '2.950.92'
$ws.Range("E46").Value = '  -2.08%  '
$ws.Range("E47").Value = '  +0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '27.42'
$ws.Range("E48").Value = '  +0.25%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '139.11'
$ws.Range("E50").Value = '  +0.75%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.46'
$ws.Range("E51").Value = '  -1.20%  '
